$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Kapil%Goyal%kapilgoyalpgi@gmail.com%1,                Poonam%Chauhan%chauhan.poonamk@gmail.com%1,                Komal%Chhikara%dimpichhikara@gmail.com%1,                Parakriti%Gupta%parakritii@gmail.com%1,                Mini P.%Singh%minipsingh@gmail.com%1]"
$ws.Range("E3").Value = "[Mohammed A.%Mamun%NULL%1,                Mark D.%Griffiths%NULL%1]"
$ws.Range("E4").Value = "[Swapnajeet%Sahoo%NULL%1,                Seema%Rani%NULL%1,                Shaheena%Parveen%NULL%1,                Ajay%Pal Singh%NULL%1,                Aseem%Mehra%NULL%1,                Subho%Chakrabarti%NULL%1,                Sandeep%Grover%NULL%1,                Cheering%Tandup%NULL%2,                Cheering%Tandup%NULL%0]"
$ws.Range("E5").Value = "[Department%of Health%coreGivesNoEmail%1,               Department%of Health of the Government of the Hong Kong Special Administrative Region%coreGivesNoEmail%1,               US%Department of Health and Human Services—Public Health Service%coreGivesNoEmail%1,               WHO%European Ministerial Conference on Mental Health%coreGivesNoEmail%1,               WHO%Website%coreGivesNoEmail%1]"
$ws.Range("E6").Value = "[Chau%PH%coreGivesNoEmail%1,               Cheung%YT%coreGivesNoEmail%1,               Yip%PSF%coreGivesNoEmail%1]"
$ws.Range("E7").Value = "[Chien-Cheng%Huang%NULL%1,                David Hung-Tsang%Yen%hjyen@vghtpe.gov.tw%1,                Hsien-Hao%Huang%NULL%1,                Wei-Fong%Kao%NULL%1,                Lee-Min%Wang%NULL%1,                Chun-I%Huang%NULL%1,                Chen-Hsen%Lee%NULL%1]"
$ws.Range("E9").Value = "[Olaoluwa%Okusaga%NULL%1,                Robert H.%Yolken%NULL%1,                Patricia%Langenberg%NULL%1,                Manana%Lapidus%NULL%1,                Timothy A.%Arling%NULL%1,                Faith B.%Dickerson%NULL%1,                Debra A.%Scrandis%NULL%1,                Emily%Severance%NULL%1,                Johanna A.%Cabassa%NULL%1,                Theodora%Balis%NULL%1,                Teodor T.%Postolache%NULL%1]"
$ws.Range("E10").Value = "[Karine%Kahil%NULL%1,                Mohamad Ali%Cheaito%NULL%1,                Rawad%El Hayek%NULL%1,                Marwa%Nofal%NULL%1,                Sarah%El Halabi%NULL%1,                Kundadak Ganesh%Kudva%NULL%1,                Victor%Pereira-Sanchez%NULL%1,                Samer%El Hayek%NULL%1]"
